# Applies the timetable room/section code corrections across the three
# timetable sheets (Regular, PreMid, PostMid) of the ECE sem5 workbook.

$wb = $excel.ActiveWorkbook

# --- Regular_Timetable ---
$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("B6").Value = "EC306 (Lab) [L206]"
$ws.Range("B7").Value = "EC306 (Lab) [L206]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PreMid_Timetable ---
$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value = "EC306 (Lab) [L206]"
$ws.Range("B9").Value = "EC306 (Lab) [L206]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PostMid_Timetable ---
$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("B8").Value = "EC306 (Lab) [L105]"
$ws.Range("B9").Value = "EC306 (Lab) [L105]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"
